$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the resolved case (old row 3, Caso 6262 - MIGUELETES 1330); all rows below shift up by one
$ws.Rows(3).Delete()

# Append the newly reported cases at the bottom of the table
$newRows = @(
    @{ row = 49; A = '7370'; B = '9/30/2025'; C = 'LA PAMPA 3621'; D = 13; E = 'Pendiente ADM'; F = 'Optical Power'; G = 'Pendiente'; H = 'Tendido a baja altura'; I = 1; J = '{"direccionesNormalizadas": [{"altura": 3621, "cod_calle": 12168, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.466333", "y": "-34.571822"}, "direccion": "LA PAMPA 3621, CABA", "nombre_calle": "LA PAMPA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'; K = -58.466333; L = -34.571822; M = 'Colegiales'; N = 'Capital Norte' },
    @{ row = 50; A = '7376'; B = '10/2/2025'; C = 'CALAZA, JOSE MARIA 1571'; D = 9; E = 'Pendiente ADM'; F = 'Optical Power'; G = 'Pendiente'; H = 'Cable cortado'; I = 1; J = '{"direccionesNormalizadas": [{"altura": 1571, "cod_calle": 3019, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.466804", "y": "-34.648650"}, "direccion": "CALAZA, JOSE MARIA 1571, CABA", "nombre_calle": "CALAZA, JOSE MARIA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'; K = -58.466804; L = -34.64865; M = 'Devoto'; N = 'Capital Norte' },
    @{ row = 51; A = '7398'; B = '10/2/2025'; C = 'GORDILLO, TIMOTEO 668'; D = 9; E = 'Pendiente ADM'; F = 'Optical Power'; G = 'Pendiente'; H = 'Tendido a baja altura'; I = 1; J = '{"direccionesNormalizadas": [{"altura": 668, "cod_calle": 7070, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.521513", "y": "-34.646130"}, "direccion": "GORDILLO, TIMOTEO 668, CABA", "nombre_calle": "GORDILLO, TIMOTEO", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'; K = -58.521513; L = -34.64613; M = 'Devoto'; N = 'Capital Norte' },
    @{ row = 52; A = '7406'; B = '10/2/2025'; C = 'FRANKLIN 871'; D = 6; E = 'Pendiente ADM'; F = 'Optical Power'; G = 'Pendiente'; H = 'Tendido a baja altura'; I = 1; J = '{"direccionesNormalizadas": [{"altura": 871, "cod_calle": 6053, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.445664", "y": "-34.609392"}, "direccion": "FRANKLIN 871, CABA", "nombre_calle": "FRANKLIN", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'; K = -58.445664; L = -34.609392; M = 'Boedo'; N = 'Capital Sur' },
    @{ row = 53; A = '7407'; B = '10/2/2025'; C = 'GAINZA, MARTIN DE, GRAL. 993'; D = 6; E = 'Pendiente ADM'; F = 'Optical Power'; G = 'Pendiente'; H = 'Tendido a baja altura'; I = 1; J = '{"direccionesNormalizadas": [{"altura": 993, "cod_calle": 7003, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.446516", "y": "-34.609820"}, "direccion": "GAINZA, MARTIN DE, GRAL. 993, CABA", "nombre_calle": "GAINZA, MARTIN DE, GRAL.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'; K = -58.446516; L = -34.60982; M = 'Boedo'; N = 'Capital Sur' },
    @{ row = 54; A = '1138'; B = '10/2/2025'; C = 'CULLEN 5085'; D = 12; E = 'Pendiente ADM'; F = 'Optical Power'; G = 'Pendiente'; H = 'Cable en panza'; I = 1; J = '{"direccionesNormalizadas": [{"altura": 5085, "cod_calle": 3202, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.487797", "y": "-34.572787"}, "direccion": "CULLEN 5085, CABA", "nombre_calle": "CULLEN", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'; K = -58.487797; L = -34.572787; M = 'Paternal'; N = 'Capital Norte' },
    @{ row = 55; A = '1146'; B = '10/2/2025'; C = 'YERBAL 1301'; D = 6; E = 'Pendiente ADM'; F = 'Optical Power'; G = 'Pendiente'; H = 'Tendido aereo cortado '; I = 1; J = '{"direccionesNormalizadas": [{"altura": 1301, "cod_calle": 26003, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.448710", "y": "-34.622159"}, "direccion": "YERBAL 1301, CABA", "nombre_calle": "YERBAL", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'; K = -58.44871; L = -34.622159; M = 'Boedo'; N = 'Capital Sur' }
)

foreach ($item in $newRows) {
    $r = $item.row
    $ws.Range("A" + $r + ":C" + $r + ",E" + $r + ":H" + $r + ",J" + $r + ":J" + $r + ",M" + $r + ":N" + $r).NumberFormat = "@"
    $ws.Range("A" + $r).Value = $item.A
    $ws.Range("B" + $r).Value = $item.B
    $ws.Range("C" + $r).Value = $item.C
    $ws.Range("D" + $r).Value = $item.D
    $ws.Range("E" + $r).Value = $item.E
    $ws.Range("F" + $r).Value = $item.F
    $ws.Range("G" + $r).Value = $item.G
    $ws.Range("H" + $r).Value = $item.H
    $ws.Range("I" + $r).Value = $item.I
    $ws.Range("J" + $r).Value = $item.J
    $ws.Range("K" + $r).Value = $item.K
    $ws.Range("L" + $r).Value = $item.L
    $ws.Range("M" + $r).Value = $item.M
    $ws.Range("N" + $r).Value = $item.N
}
